$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new characteristic text to cell E3 (shared string: "processeur : de 350 MHz à 1 GHz")
$ws.Range("E3").Value = "processeur : de 350 MHz à 1 GHz"

# Update the active selection from E2 to C3
$ws.Range("C3").Select()
